# GILDNoun.xlsx update:
#  - Adds a new shared string "Up" (sentiment up/down category)
#  - Row 3 gains two more data points: W3 (0), X3 (~0.07 momentum delta) and
#    Y3 ("Up" category)
#  - A brand new row 4 is appended with a full new day's worth of sentiment /
#    trading data

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 3 additions ----
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0.069999999999993179
$ws.Range("Y3").Value = "Up"

# ---- Row 4 (new row) ----
$ws.Range("A4").Value = 42633.884282407409
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Neutral"
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = 20775
$ws.Range("F4").Value = 3120
$ws.Range("G4").Value = 61
$ws.Range("H4").Value = 36
$ws.Range("I4").Value = 85
$ws.Range("J4").Value = 15
$ws.Range("K4").Value = 20069
$ws.Range("L4").Value = 384
$ws.Range("M4").Value = 231
$ws.Range("N4").Value = 34
$ws.Range("O4").Value = 6
$ws.Range("P4").Value = "Noun"
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = -31.95
$ws.Range("S4").Value = -0.082000000000000003
$ws.Range("S4").NumberFormat = $ws.Range("S3").NumberFormat
$ws.Range("T4").Value = -0.28000000000000003
$ws.Range("U4").Value = 6.77
$ws.Range("V4").Value = 1.88
$ws.Range("W4").Value = 0
